$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing response in row 41 (tag "เบื่อ") with the tweaked game-invite text
$ws.Range("B41").Value = "งั้นมาเล่นเกมส์ตอบคำถามกัน ตอบถูกทั้งหมด 3 ข้อและไม่ผิดเลยจะได้รางวัลจากน้องบอทแหละ <3 โอเค๊?"

# Add a new tag/response pair in row 42
$ws.Range("A42").Value = "ตกลง"
$ws.Range("B42").Value = "ข้อ 1 ประเทศอะไรรวยที่สุดดด? เพราะอะไร?"

# Update selection to match the recorded workbook view state (without
# disturbing the scrolled top-left cell of the view)
$excel.Goto($ws.Range("A43"), $false)
